$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6076
$ws1.Range("F5").Value = 366
$ws1.Range("F7").Value = 3
$ws1.Range("F9").Value = 49
$ws1.Range("F13").Value = 350
$ws1.Range("F14").Value = 607
$ws1.Range("F15").Value = 3072
$ws1.Range("F16").Value = 7
$ws1.Range("F18").Value = 1721

# Sheet "全部类型" (all types) - rows offset by +1 vs "展览"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6076
$ws4.Range("F5").Value = 366
$ws4.Range("F7").Value = 3
$ws4.Range("F10").Value = 49
$ws4.Range("F14").Value = 350
$ws4.Range("F15").Value = 607
$ws4.Range("F16").Value = 3072
$ws4.Range("F17").Value = 7
$ws4.Range("F19").Value = 1721
